$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(35, 3).Value = 1000034
$ws.Cells.Item(35, 5).Value = 'Centrale nazionale d''allarme'
$ws.Cells.Item(36, 3).Value = 1000035
$ws.Cells.Item(36, 5).Value = 'Stato maggiore federale Protezione della popolazione'
$ws.Cells.Item(37, 3).Value = 1000036
$ws.Cells.Item(37, 5).Value = 'Protezione civile'
$ws.Cells.Item(38, 3).Value = 1000038
$ws.Cells.Item(38, 5).Value = 'Comando della protezione civile'
$ws.Cells.Item(39, 3).Value = 1000039
$ws.Cells.Item(39, 5).Value = 'Impianto di protezione per la protezione della popolazione'
$ws.Cells.Item(40, 3).Value = 1000040
$ws.Cells.Item(40, 5).Value = 'Evento NBC'
$ws.Cells.Item(41, 3).Value = 1000041
$ws.Cells.Item(41, 5).Value = 'Suscettibilità da Frana'
$ws.Cells.Item(42, 3).Value = 1000042
$ws.Cells.Item(42, 5).Value = 'Pompieri'
$ws.Cells.Item(43, 3).Value = 1000043
$ws.Cells.Item(43, 5).Value = 'Vigili del fuoco'
$ws.Cells.Item(44, 3).Value = 1000044
$ws.Cells.Item(44, 5).Value = 'Struttura operativa'
$ws.Cells.Item(45, 3).Value = 1000045
$ws.Cells.Item(45, 5).Value = 'Organizzazione partner'
$ws.Cells.Item(46, 3).Value = 1000046
$ws.Cells.Item(46, 5).Value = 'Organo di condotta'
$ws.Cells.Item(47, 3).Value = 1000047
$ws.Cells.Item(47, 5).Value = 'Stato di necessità'
$ws.Cells.Item(48, 3).Value = 1000048
$ws.Cells.Item(48, 5).Value = 'Emergenza / Stato di emergenza / Evento emergenziale / Evento'
$ws.Cells.Item(49, 3).Value = 1000049
$ws.Cells.Item(49, 5).Value = 'Impianto di protezione'
$ws.Cells.Item(50, 3).Value = 1000050
$ws.Cells.Item(50, 5).Value = 'Costruzione di protezione'
$ws.Cells.Item(51, 3).Value = 1000051
$ws.Cells.Item(51, 5).Value = 'Rifugio'
$ws.Cells.Item(52, 3).Value = 1000052
$ws.Cells.Item(52, 5).Value = 'Addetto all''assistenza'
$ws.Cells.Item(53, 3).Value = 1000053
$ws.Cells.Item(53, 5).Value = 'Assistente di stato maggiore'
$ws.Cells.Item(54, 3).Value = 1000054
$ws.Cells.Item(54, 5).Value = 'Organizzazione degli Stati maggiori di condotta'
$ws.Cells.Item(55, 3).Value = 1000055
$ws.Cells.Item(55, 5).Value = 'Pioniere'
$ws.Cells.Item(56, 3).Value = 1000056
$ws.Cells.Item(56, 5).Value = 'Stato maggiore cantonale di condotta'
$ws.Cells.Item(57, 3).Value = 1000057
$ws.Cells.Item(57, 5).Value = 'Sistema d’allarme acqua'
$ws.Cells.Item(58, 3).Value = 1000058
$ws.Cells.Item(58, 5).Value = 'Allarme acqua'
$ws.Cells.Item(59, 3).Value = 1000059
$ws.Cells.Item(59, 5).Value = 'Stato maggiore regionale di condotta'
$ws.Cells.Item(60, 3).Value = 1000060
$ws.Cells.Item(60, 5).Value = 'Stato maggiore enti di primo intervento'
$ws.Cells.Item(61, 3).Value = 1000062
$ws.Cells.Item(61, 5).Value = 'Protezione della popolazione'
$ws.Cells.Item(62, 3).Value = 1000063
$ws.Cells.Item(62, 5).Value = 'Dipartimento della protezione civile'
$ws.Cells.Item(63, 3).Value = 1000064
$ws.Cells.Item(63, 5).Value = 'Legge federale sulla protezione della popolazione e sulla protezione civile'
$ws.Cells.Item(64, 3).Value = 1000066
$ws.Cells.Item(64, 5).Value = 'Volontario di protezione civile'
$ws.Cells.Item(65, 3).Value = 1000068
$ws.Cells.Item(65, 5).Value = 'Consiglio di stato'
$ws.Cells.Item(66, 3).Value = 1000069
$ws.Cells.Item(66, 5).Value = 'Milite'
$ws.Cells.Item(67, 3).Value = 1000070
$ws.Cells.Item(67, 5).Value = 'Legge sulla protezione della popolazione del 26 febbraio 2007'
$ws.Cells.Item(68, 3).Value = 1000071
$ws.Cells.Item(68, 5).Value = 'Nucleo Unitario di Valutazione e Risposta Emergenze transfrontaliere'
$ws.Cells.Item(69, 3).Value = 1000073
$ws.Cells.Item(69, 5).Value = 'Sala Operativa Regionale dell''Emergenza Urgenza'
$ws.Cells.Item(70, 3).Value = 1000074
$ws.Cells.Item(70, 5).Value = 'SOREU dei laghi'
$ws.Cells.Item(71, 3).Value = 1000075
$ws.Cells.Item(71, 5).Value = 'sezione del militare e della protezione della popolazione'
